# Apply the "disconnected_elements" diagnostic table edit.
#   B1 = 0               (bold, centered/top, thin box border)
#   A2 = 0               (bold, centered/top, thin box border)
#   B2 = "disconnected_elements" (plain, default style)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B1: value + formatting (font bold, thin border all around, centered/top) ---
$b1 = $ws.Range("B1")
$b1.Value = 0
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108   # xlCenter
$b1.VerticalAlignment = -4160     # xlTop
$b1.Borders.LineStyle = 1         # xlContinuous
$b1.Borders.Weight = 2            # xlThin

# --- A2: same formatting as B1, copied so the style is reused rather than duplicated ---
$a2 = $ws.Range("A2")
$b1.Copy()
$a2.PasteSpecial(-4122)           # xlPasteFormats
$a2.Value = 0
$excel.CutCopyMode = $false

# --- B2: plain text label (uses the shared string table, default style) ---
$ws.Range("B2").Value = "disconnected_elements"
